$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10 (Rule "R30") column C ("C1" condition value) changed from 18 to 1
$ws.Range("C10").Value = 1
